# Fruta / hortaliza, semanal
# A new weekly record was added for "Terminal Hortofrutícola Agro Chillán - Pepino dulce".
# The new record is inserted as row 26 (pushing the previous rows 26-39 down to 27-40),
# which is how the source report prepends the newest week's data above the older rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 26, shifting existing rows 26:39 down to 27:40.
$ws.Rows("26:26").Insert()

# Populate the newly inserted row 26 with this week's data.
$ws.Range("A26").Value = 7
$ws.Range("B26").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C26").Value = "Ñuble"
$ws.Range("D26").Value = 45040
$ws.Range("E26").Value = 16
$ws.Range("F26").Value = 100112043
$ws.Range("G26").Value = "Pepino dulce"
$ws.Range("H26").Value = "Cultivar IV Región"
$ws.Range("I26").Value = "Primera"
$ws.Range("J26").Value = 80
$ws.Range("K26").Value = 15000
$ws.Range("L26").Value = 16000
$ws.Range("M26").Value = 15500
$ws.Range("N26").Value = "`$/bandeja 18 kilos"
$ws.Range("O26").Value = "Provincia de Limarí"
$ws.Range("P26").Value = 861
$ws.Range("Q26").Value = 18
$ws.Range("R26").Value = "Hortaliza"
